$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Annual totals")

$ws.Range("C2").Value = 6938.57
$ws.Range("C3").Value = 31593.57
$ws.Range("C4").Value = 26095.18
$ws.Range("C5").Value = 29898.62
$ws.Range("C6").Value = 21931.64
$ws.Range("C7").Value = 36262.12
$ws.Range("C8").Value = 29402.6
$ws.Range("C9").Value = 40264.88
$ws.Range("C10").Value = 55840.55
$ws.Range("C11").Value = 22885.56
$ws.Range("C12").Value = 27443.33
$ws.Range("C13").Value = 17464.39
$ws.Range("C14").Value = 37278.68
$ws.Range("C15").Value = 10839.46
$ws.Range("C16").Value = 19073.24
$ws.Range("C17").Value = 27809.29
$ws.Range("C18").Value = 35978.57
$ws.Range("C19").Value = 28929.62
$ws.Range("C20").Value = 36575.22
$ws.Range("C21").Value = 45443.31
$ws.Range("C22").Value = 33667.9
$ws.Range("C23").Value = 4841.31
$ws.Range("C24").Value = 46983.03
$ws.Range("C25").Value = 12647.94
$ws.Range("C26").Value = 14588.06
$ws.Range("C27").Value = 24061.73
$ws.Range("C28").Value = 44763.34
$ws.Range("C29").Value = 30012.71
$ws.Range("C30").Value = 41256.33
$ws.Range("C31").Value = 26373.18
$ws.Range("C32").Value = 45724.01
$ws.Range("C33").Value = 55139.97
$ws.Range("C34").Value = 38120.45
$ws.Range("C35").Value = 37892.87
$ws.Range("C36").Value = 30379.93
$ws.Range("C37").Value = 18333.4
$ws.Range("C38").Value = 25470.42
$ws.Range("C39").Value = 8455.67
$ws.Range("C40").Value = 10870.45
$ws.Range("C41").Value = 28135.31
$ws.Range("C42").Value = 7530.32
$ws.Range("C43").Value = 20997.75
$ws.Range("C44").Value = 33761.34
$ws.Range("C45").Value = 22188.43
$ws.Range("C46").Value = 22564.18
$ws.Range("C47").Value = 29204.21
$ws.Range("C48").Value = 36147.46
$ws.Range("C49").Value = 26481.21
$ws.Range("C50").Value = 45537.84
$ws.Range("C51").Value = 34453.99
$ws.Range("C52").Value = 29874.33
$ws.Range("C53").Value = 30130.7
$ws.Range("C54").Value = 36365.53
$ws.Range("C55").Value = 40182.31
$ws.Range("C56").Value = 43866.42
$ws.Range("C57").Value = 25879
$ws.Range("C58").Value = 58971.53
$ws.Range("C59").Value = 27611.15
$ws.Range("C60").Value = 27245.76
$ws.Range("C61").Value = 26341.62
$ws.Range("C62").Value = 1791001.49